# edit.ps1
# Applies the "Major overhaul to visual input display feature" change:
#  - Renames the MCU part number PIC18F57K42 -> PIC18F47K42 (shared text)
#  - Adds a FLASH_CS / Flash Select pin entry and a UART_HOST / Write to Host
#    pin entry to the MCU-Play sheet (replacing the old, unused SS/CS row)
#  - Mirrors the UART_HOST signal and marks several rows as "O" (output) on
#    the MCU-Viz sheet
#  - Updates the remembered cell selection on both sheets

$wb = $excel.ActiveWorkbook

$wsPlay = $wb.Worksheets.Item("MCU-Play")
$wsViz  = $wb.Worksheets.Item("MCU-Viz")

# ---------------------------------------------------------------------------
# Part number rename (shared by both sheets via the F2 merged cell on
# MCU-Play; MCU-Viz's F2 uses a different, unrelated string).
# ---------------------------------------------------------------------------
if ($wsPlay.Range("F2").Value2 -eq "PIC18F57K42") {
    $wsPlay.Range("F2").Value2 = "PIC18F47K42"
}

# ---------------------------------------------------------------------------
# MCU-Play ("pinout") sheet
# ---------------------------------------------------------------------------

# Row 16: new FLASH_CS pin usage (I/J/K = Usage/IO/Notes block on the right)
$wsPlay.Range("I16").Value2 = "FLASH_CS"
$wsPlay.Range("J16").Value2 = "O"
$wsPlay.Range("K16").Value2 = "Flash Select"

# Row 17: new UART_HOST pin usage (A/B/C = Notes/IO/Usage block on the left)
$wsPlay.Range("A17").Value2 = "Write to Host"
$wsPlay.Range("B17").Value2 = "O"
$wsPlay.Range("C17").Value2 = "UART_HOST"

# Row 20: the old SS/CS usage entry on the left is removed (A/B/C cleared)
$wsPlay.Range("A20:C20").ClearContents()

# Row 20 shrinks from the default 15.75 to 15 now that it holds no wrapped
# usage text.
$wsPlay.Rows.Item(20).RowHeight = 15

# Remembered selection moves from F22 to A18
$wsPlay.Select()
$wsPlay.Range("A18").Select()

# ---------------------------------------------------------------------------
# MCU-Viz sheet
# ---------------------------------------------------------------------------

# Mark the NES/N64 strobe rows as outputs
$wsViz.Range("B5").Value2 = "O"
$wsViz.Range("B6").Value2 = "O"
$wsViz.Range("B7").Value2 = "O"
$wsViz.Range("B8").Value2 = "O"

# Row 6: mirror the new UART_HOST usage from MCU-Play (I/J/K block)
$wsViz.Range("I6").Value2 = "UART_HOST"
$wsViz.Range("J6").Value2 = "I"
$wsViz.Range("K6").Value2 = "Write to Host"

# Remembered selection moves from G21 to I19
$wsViz.Select()
$wsViz.Range("I19").Select()
